# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# - Updates the summary header values (Valor Mora, Cant. Trabajadores, Cant. Periodos)
# - Updates "Salario Basico" (column G) for the existing worker's 20 period rows
# - Re-numbers the "Periodo Mora" (column E) of those 20 rows so they run in
#   ascending chronological order (1808 .. 2003) instead of descending
# - Removes the second worker's (NIT 9005926540) 6 summary rows, which shifts
#   the signature block up from rows 46/47 to rows 40/41

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header summary values
$ws.Range("E11").Value = 624980   # VALOR MORA
$ws.Range("C13").Value = 1        # Cant. Trabajadores
$ws.Range("F13").Value = 20       # Cant. Periodos

# 2) Re-point the period column (E16:E35) to run chronologically ascending
$periods = @("1808","1809","1810","1811","1812","1901","1902","1903","1904","1905", `
             "1906","1907","1908","1909","1910","1911","1912","2001","2002","2003")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("G$row").Value = 781242
}

# 3) Drop the old NIT-9005926540 summary rows (36:41); this shifts the
#    signature rows that used to be 46/47 up to 40/41 automatically.
$ws.Rows("36:41").Delete()
